$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Full"
$ws.Range("B6").Value = $true
